$wb = $excel.ActiveWorkbook

function Set-Cell($ws, [string]$addr, $value) {
    $ws.Range($addr).Value = $value
}

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
Set-Cell $ws "H11" 435.14285
Set-Cell $ws "I11" 435.14285
Set-Cell $ws "K11" 435.14285
Set-Cell $ws "M11" -295.14285

Set-Cell $ws "H121" 1102.9375
Set-Cell $ws "J121" 1226.6923
Set-Cell $ws "L121" 3680.0769
Set-Cell $ws "N121" -7174.0769

Set-Cell $ws "H125" 1471
Set-Cell $ws "I125" 1125
Set-Cell $ws "J125" 1747.8
Set-Cell $ws "K125" 10125
Set-Cell $ws "L125" 15730.2
Set-Cell $ws "M125" -7665
Set-Cell $ws "N125" -20650.2

Set-Cell $ws "H129" 1133.475
Set-Cell $ws "I129" 1541.5
Set-Cell $ws "J129" 1061.4706
Set-Cell $ws "K129" 4624.5
Set-Cell $ws "L129" 3184.4118
Set-Cell $ws "M129" 375.5
Set-Cell $ws "N129" -13184.4118

Set-Cell $ws "H135" 597.70966
Set-Cell $ws "I135" 462.85
Set-Cell $ws "J135" 842.9091
Set-Cell $ws "K135" 4165.650000000001
Set-Cell $ws "L135" 7586.1819
Set-Cell $ws "M135" -1630.650000000001
Set-Cell $ws "N135" -12656.1819

Set-Cell $ws "H137" 2635849
Set-Cell $ws "I137" 3575928.2
Set-Cell $ws "J137" 3626.9
Set-Cell $ws "K137" 10727784.6
Set-Cell $ws "L137" 10880.7
Set-Cell $ws "M137" -10725234.6
Set-Cell $ws "N137" -15980.7

Set-Cell $ws "H138" 2900.5
Set-Cell $ws "I138" 1739
Set-Cell $ws "J138" 3948.6829
Set-Cell $ws "K138" 5217
Set-Cell $ws "L138" 11846.0487
Set-Cell $ws "M138" -77
Set-Cell $ws "N138" -22126.0487

Set-Cell $ws "H141" 461420.1
Set-Cell $ws "I141" 1166.4
Set-Cell $ws "J141" 844964.8
Set-Cell $ws "K141" 3499.2
Set-Cell $ws "L141" 2534894.4
Set-Cell $ws "M141" 1680.8
Set-Cell $ws "N141" -2545254.4

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
Set-Cell $ws "H32" 7915.53
Set-Cell $ws "I32" 7273.742
Set-Cell $ws "J32" 28666.666
Set-Cell $ws "K32" 7273.742
Set-Cell $ws "L32" 28666.666
Set-Cell $ws "M32" -6986.742
Set-Cell $ws "N32" -29240.666

Set-Cell $ws "H132" 2675.4614
Set-Cell $ws "I132" 1934.7273
Set-Cell $ws "K132" 5804.1819
Set-Cell $ws "M132" -3274.1819

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
Set-Cell $ws "H80" 634.52
Set-Cell $ws "I80" 657.7273
Set-Cell $ws "J80" 616.2857
Set-Cell $ws "K80" 657.7273
Set-Cell $ws "L80" 616.2857
Set-Cell $ws "M80" 340.2727
Set-Cell $ws "N80" -2612.2857

Set-Cell $ws "H83" 634.52
Set-Cell $ws "I83" 657.7273
Set-Cell $ws "J83" 616.2857
Set-Cell $ws "K83" 3288.6365
Set-Cell $ws "L83" 3081.4285
Set-Cell $ws "M83" 1703.3635
Set-Cell $ws "N83" -13065.4285

Set-Cell $ws "H134" 2684.83
Set-Cell $ws "I134" 2397.9607
Set-Cell $ws "J134" 10000
Set-Cell $ws "K134" 7193.882100000001
Set-Cell $ws "L134" 30000
Set-Cell $ws "M134" -4658.882100000001
Set-Cell $ws "N134" -35070

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
Set-Cell $ws "H31" 1970.2727
Set-Cell $ws "I31" 1239.3334
Set-Cell $ws "J31" 2579.389
Set-Cell $ws "K31" 1239.3334
Set-Cell $ws "L31" 2579.389
Set-Cell $ws "M31" -944.3334
Set-Cell $ws "N31" -3169.389

Set-Cell $ws "H34" 1970.2727
Set-Cell $ws "I34" 1239.3334
Set-Cell $ws "J34" 2579.389
Set-Cell $ws "K34" 1239.3334
Set-Cell $ws "L34" 2579.389
Set-Cell $ws "M34" -1037.3334
Set-Cell $ws "N34" -2983.389

Set-Cell $ws "H41" 5821.9
Set-Cell $ws "I41" 2871.4
Set-Cell $ws "J41" 8772.4
Set-Cell $ws "K41" 2871.4
Set-Cell $ws "L41" 8772.4
Set-Cell $ws "M41" -2443.4
Set-Cell $ws "N41" -9628.4

Set-Cell $ws "H50" 5424.3335
Set-Cell $ws "J50" 5424.3335
Set-Cell $ws "L50" 5424.3335
Set-Cell $ws "N50" -6674.3335

Set-Cell $ws "H51" 9515.799999999999
Set-Cell $ws "I51" 200
Set-Cell $ws "J51" 11844.75
Set-Cell $ws "K51" 200
Set-Cell $ws "L51" 11844.75
Set-Cell $ws "M51" 536
Set-Cell $ws "N51" -13316.75

Set-Cell $ws "H58" 10872985
Set-Cell $ws "I58" 2583.087
Set-Cell $ws "J58" 21743388
Set-Cell $ws "K58" 2583.087
Set-Cell $ws "L58" 21743388
Set-Cell $ws "M58" -2380.087
Set-Cell $ws "N58" -21743794

Set-Cell $ws "H59" 21701.6
Set-Cell $ws "J59" 21701.6
Set-Cell $ws "L59" 21701.6
Set-Cell $ws "N59" -23991.6

Set-Cell $ws "H60" 15781.8
Set-Cell $ws "J60" 18103
Set-Cell $ws "L60" 18103
Set-Cell $ws "N60" -19125

Set-Cell $ws "H61" 9515.799999999999
Set-Cell $ws "I61" 200
Set-Cell $ws "J61" 11844.75
Set-Cell $ws "K61" 200
Set-Cell $ws "L61" 11844.75
Set-Cell $ws "M61" 148
Set-Cell $ws "N61" -12540.75

Set-Cell $ws "H74" 16373.833
Set-Cell $ws "J74" 16373.833
Set-Cell $ws "L74" 16373.833
Set-Cell $ws "N74" -18121.833

Set-Cell $ws "H77" 16373.833
Set-Cell $ws "J77" 16373.833
Set-Cell $ws "L77" 49121.499
Set-Cell $ws "N77" -57857.499

Set-Cell $ws "H132" 3219.6667
Set-Cell $ws "I132" 2931.6
Set-Cell $ws "J132" 3699.7778
Set-Cell $ws "K132" 8794.799999999999
Set-Cell $ws "L132" 11099.3334
Set-Cell $ws "M132" -6264.799999999999
Set-Cell $ws "N132" -16159.3334

Set-Cell $ws "H134" 10206125
Set-Cell $ws "I134" 11906175
Set-Cell $ws "J134" 5828.4287
Set-Cell $ws "K134" 35718525
Set-Cell $ws "L134" 17485.2861
Set-Cell $ws "M134" -35715990
Set-Cell $ws "N134" -22555.2861

Set-Cell $ws "H136" 10872985
Set-Cell $ws "I136" 2583.087
Set-Cell $ws "J136" 21743388
Set-Cell $ws "K136" 7749.261
Set-Cell $ws "L136" 65230164
Set-Cell $ws "M136" -5199.261
Set-Cell $ws "N136" -65235264

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
Set-Cell $ws "H92" 1206.16
Set-Cell $ws "I92" 1296
Set-Cell $ws "J92" 1183.7
Set-Cell $ws "K92" 3888
Set-Cell $ws "L92" 3551.1
Set-Cell $ws "M92" -2640
Set-Cell $ws "N92" -6047.1

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
Set-Cell $ws "H102" 37279.93
Set-Cell $ws "I102" 2021.5238
Set-Cell $ws "J102" 129833.25
Set-Cell $ws "K102" 2021.5238
Set-Cell $ws "L102" 129833.25
Set-Cell $ws "M102" -399.5237999999999
Set-Cell $ws "N102" -133077.25

Set-Cell $ws "H126" 670314.25
Set-Cell $ws "I126" 2116.6667
Set-Cell $ws "J126" 1115779.4
Set-Cell $ws "K126" 6350.000100000001
Set-Cell $ws "L126" 3347338.2
Set-Cell $ws "M126" -3880.000100000001
Set-Cell $ws "N126" -3352278.2

Set-Cell $ws "H132" 3309.2979
Set-Cell $ws "I132" 3015.853
Set-Cell $ws "K132" 9047.559000000001
Set-Cell $ws "M132" -6517.559000000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
Set-Cell $ws "H61" 90913700
Set-Cell $ws "I61" 250002700
Set-Cell $ws "J61" 5697.5713
Set-Cell $ws "K61" 250002700
Set-Cell $ws "L61" 5697.5713
Set-Cell $ws "M61" -250002498
Set-Cell $ws "N61" -6101.5713

Set-Cell $ws "H113" 90913700
Set-Cell $ws "I113" 250002700
Set-Cell $ws "J113" 5697.5713
Set-Cell $ws "K113" 250002700
Set-Cell $ws "L113" 5697.5713
Set-Cell $ws "M113" -250000530
Set-Cell $ws "N113" -10037.5713

Set-Cell $ws "H132" 2886.5625
Set-Cell $ws "I132" 1947.1428
Set-Cell $ws "J132" 4680
Set-Cell $ws "K132" 5841.428400000001
Set-Cell $ws "L132" 14040
Set-Cell $ws "M132" -3311.428400000001
Set-Cell $ws "N132" -19100

Set-Cell $ws "H136" 3452746.5
Set-Cell $ws "I136" 5559230.5
Set-Cell $ws "J136" 5773.1816
Set-Cell $ws "K136" 16677691.5
Set-Cell $ws "L136" 17319.5448
Set-Cell $ws "M136" -16675141.5
Set-Cell $ws "N136" -22419.5448

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
Set-Cell $ws "H113" 2306
Set-Cell $ws "I113" 228.4
Set-Cell $ws "K113" 685.2
Set-Cell $ws "M113" 1484.8

Set-Cell $ws "H132" 238982.69
Set-Cell $ws "I132" 314822.97
Set-Cell $ws "J132" 36741.918
Set-Cell $ws "K132" 944468.9099999999
Set-Cell $ws "L132" 110225.754
Set-Cell $ws "M132" -941938.9099999999
Set-Cell $ws "N132" -115285.754
